# "Generate Report for Handback"
#
# The localization-status CI report is regenerated: the two files that were
# still "Ready for handoff" (5696aa25-...md and c9f01839-...md) have now
# been handed back (for both the zh-cn and de-de target languages), so:
#   - their Status flips to "Handed back: in sync with en-US"
#   - their Latest Target File / Latest Handback File columns get populated
#     (with a hyperlink on the target file, same as the existing Source
#     File Name hyperlinks)
#   - their Latest Handback DateTime changes from the placeholder
#     "0001-01-01 00:00:00" to the real handback timestamp

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: just the rolled-up status text for the two files.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E4").Value = $statusHandedBack
$overview.Range("F4").Value = $statusHandedBack
$overview.Range("E5").Value = $statusHandedBack
$overview.Range("F5").Value = $statusHandedBack

# ---------------------------------------------------------------------
# zh-cn sheet: rows 4 (5696aa25...) and 5 (c9f01839...)
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C4").Value = $statusHandedBack
$zhcn.Range("J4").Value = $zhcn.Range("G4").Value()
$zhcn.Range("K4").Value = "2016-09-03 00:28:32"
$zhcn.Range("I4").Value = "5696aa25-bf1e-463d-b662-d1da1c4e31e4.md"
$zhcn.Hyperlinks.Add($zhcn.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/e0cfb22610d478f6d07681049d1ca1a575d946bf/e2e/5696aa25-bf1e-463d-b662-d1da1c4e31e4.md", "", "", "5696aa25-bf1e-463d-b662-d1da1c4e31e4.md")

$zhcn.Range("C5").Value = $statusHandedBack
$zhcn.Range("J5").Value = $zhcn.Range("G5").Value()
$zhcn.Range("K5").Value = "2016-09-03 00:28:32"
$zhcn.Range("I5").Value = "c9f01839-a57e-4399-84fe-18a65f9ef9b3.md"
$zhcn.Hyperlinks.Add($zhcn.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/e0cfb22610d478f6d07681049d1ca1a575d946bf/e2e/c9f01839-a57e-4399-84fe-18a65f9ef9b3.md", "", "", "c9f01839-a57e-4399-84fe-18a65f9ef9b3.md")

# ---------------------------------------------------------------------
# de-de sheet: rows 4 (5696aa25...) and 5 (c9f01839...)
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C4").Value = $statusHandedBack
$dede.Range("J4").Value = $dede.Range("G4").Value()
$dede.Range("K4").Value = "2016-09-03 00:28:39"
$dede.Range("I4").Value = "5696aa25-bf1e-463d-b662-d1da1c4e31e4.md"
$dede.Hyperlinks.Add($dede.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/0b3663163da692da6758182ea8b0b3878e426b4a/e2e/5696aa25-bf1e-463d-b662-d1da1c4e31e4.md", "", "", "5696aa25-bf1e-463d-b662-d1da1c4e31e4.md")

$dede.Range("C5").Value = $statusHandedBack
$dede.Range("J5").Value = $dede.Range("G5").Value()
$dede.Range("K5").Value = "2016-09-03 00:28:39"
$dede.Range("I5").Value = "c9f01839-a57e-4399-84fe-18a65f9ef9b3.md"
$dede.Hyperlinks.Add($dede.Range("I5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/0b3663163da692da6758182ea8b0b3878e426b4a/e2e/c9f01839-a57e-4399-84fe-18a65f9ef9b3.md", "", "", "c9f01839-a57e-4399-84fe-18a65f9ef9b3.md")
